# Auto-generated edit script: update cryptos Price (D) and Volume(1h) (E) columns
# to reflect the latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.489.73"
$ws.Range("E2").Value = "  -6.79%  "
$ws.Range("D3").Value = "2.425.24"
$ws.Range("E3").Value = "  -10.90%  "
$ws.Range("E4").Value = "  +0.01%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "466.41"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -7.37%  "
$ws.Range("E6").Value = "  -5.18%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  -6.85%  "
$ws.Range("D9").Value = "2.441.57"
$ws.Range("E9").Value = "  -10.69%  "
$ws.Range("E10").Value = "  -8.96%  "
$ws.Range("E11").Value = "  -12.04%  "
$ws.Range("E12").Value = "  -9.12%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.121"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -3.95%  "
$ws.Range("D14").Value = "2.845.74"
$ws.Range("E14").Value = "  -11.18%  "
$ws.Range("D15").Value = "54.442.69"
$ws.Range("E15").Value = "  -7.09%  "
$ws.Range("E16").Value = "  -1.34%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "19.79"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -8.58%  "
$ws.Range("D18").Value = "2.442.52"
$ws.Range("E18").Value = "  -10.67%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "4.20"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -11.80%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "310.30"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -9.57%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "9.48"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -13.40%  "
$ws.Range("E22").Value = "  +0.17%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.66"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.91%  "
$ws.Range("E24").Value = "  -13.97%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "56.26"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -10.52%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.01"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("E27").Value = "  -9.64%  "
$ws.Range("E28").Value = "  -10.01%  "
$ws.Range("D29").Value = "2.527.28"
$ws.Range("E29").Value = "  -11.20%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "7.12"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -5.19%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  -13.16%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "147.17"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -3.19%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "17.76"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -7.14%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.44"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -10.15%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "4.99"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -7.94%  "
$ws.Range("E37").Value = "  -15.27%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.05"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -6.52%  "
$ws.Range("E39").Value = "  -14.77%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.993"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.29%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "32.90"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -8.21%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.595"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("E43").Value = "  -6.10%  "
$ws.Range("E44").Value = "  -8.48%  "
$ws.Range("E45").Value = "  -10.90%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "10.07"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("D47").Value = "1.935.48"
$ws.Range("E47").Value = "  -11.70%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.0883"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.22%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.0217"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -4.02%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "4.21"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -11.21%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "16.59"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -12.68%  "
